# Update the "Глобальные показатели ЦУР" sheet: add a new 2020 data column (Q)
# mirroring the existing 2019 column (P) formatting, and move the active
# selection to N19, matching the authored workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write the new 2020 values into column Q -----------------------------
$ws.Range("Q4").Value  = 2020
$ws.Range("Q5").Value  = 0.02
$ws.Range("Q6").Value  = 0
$ws.Range("Q7").Value  = 0
$ws.Range("Q8").Value  = 0
$ws.Range("Q9").Value  = 0.54
$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("Q14").Value = 0

# --- 2. Mirror the formatting of column P (2019) onto the new column Q -----
$ws.Range("P4:P14").Copy() | Out-Null
$ws.Range("Q4:Q14").PasteSpecial(-4122) | Out-Null

# --- 3. Restore the active cell selection to N19 ----------------------------
$ws.Range("N19").Select() | Out-Null
